$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.055906665083465
$ws.Range("D2").Value = 1.06220660608802
$ws.Range("E2").Value = 1.069091184673903
$ws.Range("F2").Value = 1.075434052130314
$ws.Range("I2").Value = 1.048407290474736
$ws.Range("J2").Value = 1.060911076716469
$ws.Range("K2").Value = 1.064928581078178
$ws.Range("L2").Value = 1.071794619448596
$ws.Range("M2").Value = 1.078120632600537
$ws.Range("N2").Value = 1.062417692133608

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.056869215008701
$ws.Range("D3").Value = 1.062984887817917
$ws.Range("E3").Value = 1.070016515963961
$ws.Range("F3").Value = 1.076380801812208
$ws.Range("I3").Value = 1.048659195798006
$ws.Range("J3").Value = 1.061525427898315
$ws.Range("K3").Value = 1.065521725564153
$ws.Range("L3").Value = 1.072535783175583
$ws.Range("M3").Value = 1.078884378970342
$ws.Range("N3").Value = 1.0630329157646

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.057492631994108
$ws.Range("D4").Value = 1.063488984249987
$ws.Range("E4").Value = 1.070616182178601
$ws.Range("F4").Value = 1.076994362476733
$ws.Range("I4").Value = 1.048821296296379
$ws.Range("J4").Value = 1.061922889447323
$ws.Range("K4").Value = 1.065905363666561
$ws.Range("L4").Value = 1.073015646831872
$ws.Range("M4").Value = 1.07937889358811
$ws.Range("N4").Value = 1.063430941754596

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.05775485501329
$ws.Range("D5").Value = 1.06370102382012
$ws.Range("E5").Value = 1.070868499557184
$ws.Range("F5").Value = 1.077252529501843
$ws.Range("I5").Value = 1.048889227566463
$ws.Range("J5").Value = 1.062089965824141
$ws.Range("K5").Value = 1.06606660437017
$ws.Range("L5").Value = 1.073217447695625
$ws.Range("M5").Value = 1.079586862879269
$ws.Range("N5").Value = 1.063598255399034

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.057798891505615
$ws.Range("D6").Value = 1.063736633044203
$ws.Range("E6").Value = 1.070910877482269
$ws.Range("F6").Value = 1.077295890101033
$ws.Range("I6").Value = 1.048900620856706
$ws.Range("J6").Value = 1.06211801769006
$ws.Range("K6").Value = 1.066093674974799
$ws.Range("L6").Value = 1.073251334796844
$ws.Range("M6").Value = 1.079621786222491
$ws.Range("N6").Value = 1.063626347101819

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.057496135287891
$ws.Range("D7").Value = 1.063491817073358
$ws.Range("E7").Value = 1.070619552802028
$ws.Range("F7").Value = 1.076997811230384
$ws.Range("I7").Value = 1.048822204845864
$ws.Range("J7").Value = 1.061925121995912
$ws.Range("K7").Value = 1.065907518334124
$ws.Range("L7").Value = 1.073018343046246
$ws.Range("M7").Value = 1.07938167218937
$ws.Range("N7").Value = 1.06343317747366

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.056231842051894
$ws.Range("D8").Value = 1.062469526513669
$ws.Range("E8").Value = 1.069403714447809
$ws.Range("F8").Value = 1.075753813015821
$ws.Range("I8").Value = 1.048492608898031
$ws.Range("J8").Value = 1.061118712404145
$ws.Range("K8").Value = 1.065129070810348
$ws.Range("L8").Value = 1.072045040303831
$ws.Range("M8").Value = 1.078378677368114
$ws.Range("N8").Value = 1.062625622687777

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.054008504217066
$ws.Range("D9").Value = 1.060671978749216
$ws.Range("E9").Value = 1.067268319356703
$ws.Range("F9").Value = 1.073569068719205
$ws.Range("I9").Value = 1.047904956487407
$ws.Range("J9").Value = 1.059697262146265
$ws.Range("K9").Value = 1.063756119502309
$ws.Range("L9").Value = 1.070332159802824
$ws.Range("M9").Value = 1.076613772536863
$ws.Range("N9").Value = 1.061202153807487

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.052529369268075
$ws.Range("D10").Value = 1.059476287282398
$ws.Range("E10").Value = 1.065849554084721
$ws.Range("F10").Value = 1.072117589015011
$ws.Range("I10").Value = 1.047508606095251
$ws.Range("J10").Value = 1.058749380195885
$ws.Range("K10").Value = 1.062840052721853
$ws.Range("L10").Value = 1.069191787441591
$ws.Range("M10").Value = 1.075438919860004
$ws.Range("N10").Value = 1.060252925756014

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.051889632558652
$ws.Range("D11").Value = 1.058959190379322
$ws.Range("E11").Value = 1.065236375590872
$ws.Range("F11").Value = 1.071490289144945
$ws.Range("I11").Value = 1.0473359017399
$ws.Range("J11").Value = 1.058338889619042
$ws.Range("K11").Value = 1.062443216336291
$ws.Range("L11").Value = 1.068698374601739
$ws.Range("M11").Value = 1.074930625333327
$ws.Range("N11").Value = 1.059841852235469

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.051652117734044
$ws.Range("D12").Value = 1.058767215839881
$ws.Range("E12").Value = 1.065008788817869
$ws.Range("F12").Value = 1.071257463738909
$ws.Range("I12").Value = 1.047271589591383
$ws.Range("J12").Value = 1.058186408325289
$ws.Range("K12").Value = 1.062295788540238
$ws.Range("L12").Value = 1.06851515666548
$ws.Range("M12").Value = 1.074741887015865
$ws.Range("N12").Value = 1.059689154400788

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.051703060392761
$ws.Range("D13").Value = 1.058808390562088
$ws.Range("E13").Value = 1.065057599012081
$ws.Range("F13").Value = 1.07130739733506
$ws.Range("I13").Value = 1.047285392091107
$ws.Range("J13").Value = 1.058219116374202
$ws.Range("K13").Value = 1.062327413423003
$ws.Range("L13").Value = 1.06855445491253
$ws.Range("M13").Value = 1.074782369070463
$ws.Range("N13").Value = 1.059721908898883

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.051869997221598
$ws.Range("D14").Value = 1.058943319691655
$ws.Range("E14").Value = 1.065217559619245
$ws.Range("F14").Value = 1.071471040009904
$ws.Range("I14").Value = 1.047330588980739
$ws.Range("J14").Value = 1.058326285601504
$ws.Range("K14").Value = 1.062431030415002
$ws.Range("L14").Value = 1.068683228569449
$ws.Range("M14").Value = 1.074915022834581
$ws.Range("N14").Value = 1.05982923031878

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.05197286736356
$ws.Range("D15").Value = 1.059026467043121
$ws.Range("E15").Value = 1.065316139859974
$ws.Range("F15").Value = 1.071571889779614
$ws.Range("I15").Value = 1.047358414816446
$ws.Range("J15").Value = 1.05839231522355
$ws.Range("K15").Value = 1.062494868967274
$ws.Range("L15").Value = 1.06876257792951
$ws.Range("M15").Value = 1.074996763830869
$ws.Range("N15").Value = 1.059895353710463

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.052571842124351
$ws.Range("D16").Value = 1.059510619010404
$ws.Range("E16").Value = 1.065890273194184
$ws.Range("F16").Value = 1.072159246247782
$ws.Range("I16").Value = 1.047520045160215
$ws.Range("J16").Value = 1.05877662209756
$ws.Range("K16").Value = 1.062866385873991
$ws.Range("L16").Value = 1.069224541639397
$ws.Range("M16").Value = 1.075472662732886
$ws.Range("N16").Value = 1.060280206344315

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.052947761531088
$ws.Range("D17").Value = 1.059814488463453
$ws.Range("E17").Value = 1.066250722346723
$ws.Range("F17").Value = 1.072528001778868
$ws.Range("I17").Value = 1.047621142269153
$ws.Range("J17").Value = 1.059017674517692
$ws.Range("K17").Value = 1.063099382861695
$ws.Range("L17").Value = 1.069514420678104
$ws.Range("M17").Value = 1.075771295974728
$ws.Range("N17").Value = 1.060521601086533

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.053167100151665
$ws.Range("D18").Value = 1.059991792557363
$ws.Range("E18").Value = 1.066461077710599
$ws.Range("F18").Value = 1.072743206373978
$ws.Range("I18").Value = 1.047680006069424
$ws.Range("J18").Value = 1.059158271291841
$ws.Range("K18").Value = 1.063235269265035
$ws.Range("L18").Value = 1.069683538414769
$ws.Range("M18").Value = 1.075945524567797
$ws.Range("N18").Value = 1.060662397524227

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.053241901000201
$ws.Range("D19").Value = 1.060052259178884
$ws.Range("E19").Value = 1.066532822321278
$ws.Range("F19").Value = 1.072816605168867
$ws.Range("I19").Value = 1.047700059375335
$ws.Range("J19").Value = 1.059206210278597
$ws.Range("K19").Value = 1.063281600138572
$ws.Range("L19").Value = 1.069741209299403
$ws.Range("M19").Value = 1.076004938914285
$ws.Range("N19").Value = 1.060710404589843

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.052907421566117
$ws.Range("D20").Value = 1.05978187970557
$ws.Range("E20").Value = 1.066212038024408
$ws.Range("F20").Value = 1.072488425819734
$ws.Range("I20").Value = 1.047610306304511
$ws.Range("J20").Value = 1.058991812387391
$ws.Range("K20").Value = 1.06307438622047
$ws.Range("L20").Value = 1.069483315648687
$ws.Range("M20").Value = 1.075739251217949
$ws.Range("N20").Value = 1.060495702229039

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.051820835390582
$ws.Range("D21").Value = 1.058903583702635
$ws.Range("E21").Value = 1.065170450365362
$ws.Range("F21").Value = 1.071422846294315
$ws.Range("I21").Value = 1.047317284092899
$ws.Range("J21").Value = 1.058294727116019
$ws.Range("K21").Value = 1.062400518478406
$ws.Range("L21").Value = 1.06864530633938
$ws.Range("M21").Value = 1.074875957806626
$ws.Range("N21").Value = 1.059797627016626

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.051138303263864
$ws.Range("D22").Value = 1.05835193337984
$ws.Range("E22").Value = 1.064516576134291
$ws.Range("F22").Value = 1.070753926061116
$ws.Range("I22").Value = 1.04713211201168
$ws.Range("J22").Value = 1.057856402954333
$ws.Range("K22").Value = 1.061976686379125
$ws.Range("L22").Value = 1.068118750199124
$ws.Range("M22").Value = 1.074333547039346
$ws.Range("N22").Value = 1.059358680384354

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.051500064750568
$ws.Range("D23").Value = 1.058644319258562
$ws.Range("E23").Value = 1.064863110751138
$ws.Range("F23").Value = 1.071108433174037
$ws.Range("I23").Value = 1.047230363921155
$ws.Range("J23").Value = 1.058088770275954
$ws.Range("K23").Value = 1.062201381188551
$ws.Range("L23").Value = 1.068397855590148
$ws.Range("M23").Value = 1.07462105322351
$ws.Range("N23").Value = 1.059591377694224

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.052925649239964
$ws.Range("D24").Value = 1.059796614007631
$ws.Range("E24").Value = 1.066229517460216
$ws.Range("F24").Value = 1.072506308134678
$ws.Range("I24").Value = 1.047615202933227
$ws.Range("J24").Value = 1.059003498386456
$ws.Range("K24").Value = 1.063085681178663
$ws.Range("L24").Value = 1.069497370560835
$ws.Range("M24").Value = 1.075753730737996
$ws.Range("N24").Value = 1.060507404823564

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.054582750007111
$ws.Range("D25").Value = 1.061136222481465
$ws.Range("E25").Value = 1.067819524441912
$ws.Range("F25").Value = 1.074132999298033
$ws.Range("I25").Value = 1.048057688680525
$ws.Range("J25").Value = 1.060064789045296
$ws.Range("K25").Value = 1.064111199669926
$ws.Range("L25").Value = 1.070774712515978
$ws.Range("M25").Value = 1.077069739297706
$ws.Range("N25").Value = 1.061570202636871

